# Generate Report for Handoff
#
# This script updates the localization-status workbook to reflect a new
# handoff cycle:
#   - the source file guid is renamed from 0061f551-...-c4403e965861.md
#     to 617b5400-...-0532877aa549.md
#   - a brand-new dependency file ffff30bae099-...-fa1bd61ee460.md is added,
#     which piggy-backs on the same (renamed) translated .xlf files
#   - the .localization-config row moves down one row on every sheet
#   - the Overview sheet gains a third data row, and the per-locale sheets
#     (zh-cn / de-de) gain a fourth row, shifting the dimension accordingly

$wb = $excel.ActiveWorkbook

# Source guid was renamed as part of this handoff cycle; a second guid
# (depGuid) is a brand-new file that depends on / reuses the same
# translated artifacts as the renamed source file.
$newGuid = "617b5400-e27b-4269-bf2c-0532877aa549"
$depGuid = "ffff30bae099-e70f-4860-adae-fa1bd61ee460"

$newHash = "3485b3f130250f9670cc4c318f4fa47dd7fafe79"

$zhXlf = "$newGuid.$newHash.zh-cn.xlf"
$deXlf = "$newGuid.$newHash.de-de.xlf"

$handoffDtZh = "2016-03-04 11:06:43"
$handoffDtDe = "2016-03-04 11:07:01"
$neverDt = "0001-01-01 00:00:00"

$commitRepo = "0a93c757523dfb6d4fc27fe3ba21f5c726cf6e17"
$commitZh = "ca17233735339a4538115f6fa125f98a5267ebfe"
$commitDe = "03a07ddeab8a3576f80a33a5e4616d84b240ab17"

function Set-CellLink {
    param($ws, [string]$cellRef, [string]$text, [string]$url)
    $ws.Range($cellRef).Value = $text
    [void]$ws.Hyperlinks.Add($ws.Range($cellRef), $url, $null, $null, $text)
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

Set-CellLink $ws1 "A2" "$newGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$newGuid.md"
Set-CellLink $ws1 "A3" "$depGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$depGuid.md"
Set-CellLink $ws1 "A4" ".localization-config" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/.localization-config"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("A1").Value = "Source File Name"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Latest Handoff File"
$ws2.Range("D1").Value = "Latest Handoff Datetime"
$ws2.Range("E1").Value = "Latest Target File"
$ws2.Range("F1").Value = "Latest Handback File"
$ws2.Range("G1").Value = "Latest Handback DateTime"
$ws2.Range("H1").Value = "Handoff Reason"
$ws2.Range("I1").Value = "Dependency From"

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = $handoffDtZh
$ws2.Range("G2").Value = $neverDt
$ws2.Range("H2").Value = "Include"

$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("D3").Value = $handoffDtZh
$ws2.Range("G3").Value = $neverDt
$ws2.Range("H3").Value = "Include"

$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = $neverDt
$ws2.Range("G4").Value = $neverDt
$ws2.Range("H4").Value = "Ignored"

Set-CellLink $ws2 "A2" "$newGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$newGuid.md"
Set-CellLink $ws2 "C2" "$zhXlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/$zhXlf"
Set-CellLink $ws2 "A3" "$depGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$depGuid.md"
Set-CellLink $ws2 "C3" "$zhXlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/$zhXlf"
Set-CellLink $ws2 "A4" ".localization-config" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/.localization-config"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("A1").Value = "Source File Name"
$ws3.Range("B1").Value = "Status"
$ws3.Range("C1").Value = "Latest Handoff File"
$ws3.Range("D1").Value = "Latest Handoff Datetime"
$ws3.Range("E1").Value = "Latest Target File"
$ws3.Range("F1").Value = "Latest Handback File"
$ws3.Range("G1").Value = "Latest Handback DateTime"
$ws3.Range("H1").Value = "Handoff Reason"
$ws3.Range("I1").Value = "Dependency From"

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = $handoffDtDe
$ws3.Range("G2").Value = $neverDt
$ws3.Range("H2").Value = "Include"

$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("D3").Value = $handoffDtDe
$ws3.Range("G3").Value = $neverDt
$ws3.Range("H3").Value = "Include"

$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = $neverDt
$ws3.Range("G4").Value = $neverDt
$ws3.Range("H4").Value = "Ignored"

Set-CellLink $ws3 "A2" "$newGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$newGuid.md"
Set-CellLink $ws3 "C2" "$deXlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/$deXlf"
Set-CellLink $ws3 "A3" "$depGuid.md" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/e2e/$depGuid.md"
Set-CellLink $ws3 "C3" "$deXlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/$deXlf"
Set-CellLink $ws3 "A4" ".localization-config" "https://github.com/OpenLocalizationTest/oltest/blob/$commitRepo/.localization-config"

Write-Host "Report generated for handoff."
